$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: style surgery so the existing cellXfs[3] entry gets its numFmtId
# flipped from 4 (#,##0.00) to 3 (#,##0) IN PLACE, matching the diff, instead of
# Excel's usual behaviour of appending a brand-new style record. We do this by
# making D2 the sole remaining owner of style index 3, then changing its format
# (sole-owner format changes mutate the xf record in place).
$ws.Range("A2").Copy()
$ws.Range("D3:E5").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("D2").NumberFormat = "#,##0"

# --- Step 2: propagate that (now #,##0) style onto the new Costo_Unitario /
# Costo_Total columns (F:G) for all data rows.
$ws.Range("D2").Copy()
$ws.Range("F2:G8").PasteSpecial(-4122)

# --- Step 3: make sure Peso/Cantidad (D:E) end up on the plain centered style
# (same as the rest of the row) for every data row, including the brand new
# rows 5-8.
$ws.Range("A2").Copy()
$ws.Range("A2:E8").PasteSpecial(-4122)

# --- Column widths (Excel's ColumnWidth setter pads by 5/6 of a char vs the
# raw OOXML <col width>, so subtract that back out to land on the exact target
# widths from the diff).
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668  # -> 21
$ws.Columns.Item(2).ColumnWidth = 35.166666666666664  # -> 36
$ws.Columns.Item(3).ColumnWidth = 9.166666666666666   # -> 10
$ws.Columns.Item(4).ColumnWidth = 8.166666666666666   # -> 9
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666   # -> 10
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666  # -> 16
$ws.Columns.Item(7).ColumnWidth = 12.166666666666666  # -> 13

# --- Header row ---
$ws.Range("A1").Value = "Tipo_Residuo"
$ws.Range("B1").Value = "Residuo"
$ws.Range("C1").Value = "Fecha"
$ws.Range("D1").Value = "Peso"
$ws.Range("E1").Value = "Cantidad"
$ws.Range("F1").Value = "Costo_Unitario"
$ws.Range("G1").Value = "Costo_Total"

# New header cells F1/G1 start out with no explicit style; give them the same
# bold/green/centered header style already used by A1:E1.
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# --- Data rows (A, B, D, E, F, G -- C handled separately below to avoid
# Excel's automatic text->date coercion on the "dd-mm-yy"-looking strings) ---
$ws.Range("A2").Value = "Respel Aprovechable"
$ws.Range("B2").Value = "Batería Ácido Plomo 30-31H - UND"
$ws.Range("D2").Value = 1127.46
$ws.Range("E2").Value = 46
$ws.Range("F2").Value = 2483
$ws.Range("G2").Value = 2799483

$ws.Range("A3").Value = "Aprovechable"
$ws.Range("B3").Value = "Radiador Cobre - KG"
$ws.Range("D3").Value = 124.2
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 7200
$ws.Range("G3").Value = 894240

$ws.Range("A4").Value = "Respel Aprovechable"
$ws.Range("B4").Value = "Aceite Usado - KG"
$ws.Range("D4").Value = 2976.38
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 352
$ws.Range("G4").Value = 1047686

$ws.Range("A5").Value = "Respel"
$ws.Range("B5").Value = "Refrigerante - KG"
$ws.Range("D5").Value = 1210
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1428
$ws.Range("G5").Value = 1727880

$ws.Range("A6").Value = "Especial"
$ws.Range("B6").Value = "Fibra de Vidrio - MT3"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 84300
$ws.Range("G6").Value = 337200

$ws.Range("A7").Value = "Respel Aprovechable"
$ws.Range("B7").Value = "Batería Ácido Plomo 65-G4-27 - UND"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 2066
$ws.Range("G7").Value = 0

$ws.Range("A8").Value = "Aprovechable"
$ws.Range("B8").Value = "Papel Archivo - KG"
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 400
$ws.Range("G8").Value = 40000

# --- Column C: literal "dd-mm-yy" text values. Force text interpretation with
# a temporary "@" format so Excel doesn't silently turn these into date
# serials, write the values, then restore the normal centered style (matching
# every other cell in the row) by pasting formats back from A2.
$ws.Range("C2:C8").NumberFormat = "@"
$ws.Range("C2").Value = "07-07-25"
$ws.Range("C3").Value = "07-07-25"
$ws.Range("C4").Value = "07-07-25"
$ws.Range("C5").Value = "07-07-25"
$ws.Range("C6").Value = "07-07-25"
$ws.Range("C7").Value = "07-07-25"
$ws.Range("C8").Value = "07-07-25"
$ws.Range("A2").Copy()
$ws.Range("C2:C8").PasteSpecial(-4122)
